# Workbook under edit: VSP-Data-Prod.xlsx
# Update the "last run" timestamp values recorded in column B of each
# Katalon bootstrap sheet, reflecting a fresh test execution.
$wb = $excel.ActiveWorkbook

$wsACH = $wb.Worksheets.Item("AddModifyDeleteACH")
$wsACH.Range("B2").Value = "Thu Sep 04 06:35:53 IST 2025"
$wsACH.Range("B3").Value = "Thu Sep 04 06:37:10 IST 2025"
$wsACH.Range("B4").Value = "Thu Sep 04 06:38:18 IST 2025"

$wsCC = $wb.Worksheets.Item("AddModifyDeleteCC")
$wsCC.Range("B2").Value = "Thu Sep 04 06:39:35 IST 2025"

$wsProfile = $wb.Worksheets.Item("CreateModifyDeleteProfile")
$wsProfile.Range("B2").Value = "Thu Sep 04 06:40:50 IST 2025"
